# The workbook has a single sheet of daily price records (grapes / "Uva")
# for "Vega Monumental Concepción". A new daily record was inserted as a
# new row 140 (pushing the former rows 140-193 down to 141-194), so the
# worksheet dimension grows from A1:T193 to A1:T194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 140, shifting rows 140-193 down to 141-194.
$ws.Rows.Item(140).Insert(1)

# Fill in the values for the newly inserted row 140.
$ws.Range("A140").Value = 11
$ws.Range("B140").Value = "Vega Monumental Concepción"
$ws.Range("C140").Value = "Bíobío"
$ws.Range("D140").Value2 = 44992
$ws.Range("E140").Value = 8
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100109
$ws.Range("H140").Value = "Uva"
$ws.Range("I140").Value = 100109001
$ws.Range("J140").Value = "Uva"
$ws.Range("K140").Value = "Superior Seedless"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 250
$ws.Range("N140").Value = 10000
$ws.Range("O140").Value = 11000
$ws.Range("P140").Value = 10400
$ws.Range("Q140").Value = "$/bandeja 18 kilos"
$ws.Range("R140").Value = "Región de O'Higgins"
$ws.Range("S140").Value = 578
$ws.Range("T140").Value = 18
